$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates derived from the new TPM (transcripts-per-million) values.
# Ligand (Sema3d) columns G,H depend on the sending cluster; Receptor (Nrp1)
# columns M,N depend on the target cluster. Columns I,J,O,P,Q,R,S,T are all
# derived (specificity / products) from G,H,M,N, so every dependent cell is
# updated here to keep the sheet internally consistent with the diff.
$updates = @{
    "G2" = 0.01606133333333333
    "H2" = 0.048184
    "I2" = 0.009938633312098436
    "J2" = 0.0127707462341618
    "M2" = 63.46725166666666
    "N2" = 190.401755
    "O2" = 0.2354497988808272
    "P2" = 0.2397164477183668
    "Q2" = 1.019368684768889
    "R2" = 9.174318162919999
    "S2" = 0.002340049214483867
    "T2" = 0.003061357921965976
    "G3" = 0.01606133333333333
    "H3" = 0.048184
    "I3" = 0.009938633312098436
    "J3" = 0.0127707462341618
    "O3" = 0.1779985000094065
    "P3" = 0.1812240584798697
    "Q3" = 0.7706360239333334
    "R3" = 6.9357242154
    "S3" = 0.001769061821697041
    "T3" = 0.002314366462371314
    "G4" = 0.01606133333333333
    "H4" = 0.048184
    "I4" = 0.009938633312098436
    "J4" = 0.0127707462341618
    "M4" = 64.53809233333334
    "N4" = 193.614277
    "O4" = 0.2394223865221556
    "P4" = 0.243761023683841
    "Q4" = 1.036567813663111
    "R4" = 9.329110322968001
    "S4" = 0.002379531306351204
    "T4" = 0.003113010175245838
    "G5" = 0.01606133333333333
    "H5" = 0.048184
    "I5" = 0.009938633312098436
    "J5" = 0.0127707462341618
    "M5" = 14.3933435
    "N5" = 28.786687
    "O5" = 0.0533961963580272
    "P5" = 0.03624253541791403
    "Q5" = 0.2311762877346667
    "R5" = 1.387057726408
    "S5" = 0.0005306852158632383
    "T5" = 0.0004628442227048012
    "G6" = 0.01606133333333333
    "H6" = 0.048184
    "I6" = 0.009938633312098436
    "J6" = 0.0127707462341618
    "M6" = 79.17795566666666
    "N6" = 237.533867
    "O6" = 0.2937331182295834
    "P6" = 0.2990559347000084
    "Q6" = 1.271703538614222
    "R6" = 11.445331847528
    "S6" = 0.002919305753703086
    "T6" = 0.003819167451873869
    "I7" = 0.3247643973172043
    "J7" = 0.4173092591090693
    "M7" = 63.46725166666666
    "N7" = 190.401755
    "O7" = 0.2354497988808272
    "P7" = 0.2397164477183668
    "Q7" = 33.30987733997611
    "R7" = 299.788896059785
    "S7" = 0.07646571203198882
    "T7" = 0.1000358931936096
    "I8" = 0.3247643973172043
    "J8" = 0.4173092591090693
    "O8" = 0.1779985000094065
    "P8" = 0.1812240584798697
    "S8" = 0.05780757557892127
    "T8" = 0.07562647757697306
    "I9" = 0.3247643973172043
    "J9" = 0.4173092591090693
    "M9" = 64.53809233333334
    "N9" = 193.614277
    "O9" = 0.2394223865221556
    "P9" = 0.243761023683841
    "Q9" = 33.87189271515989
    "R9" = 304.8470344364391
    "S9" = 0.07775586706311462
    "T9" = 0.101723732193172
    "I10" = 0.3247643973172043
    "J10" = 0.4173092591090693
    "M10" = 14.3933435
    "N10" = 28.786687
    "O10" = 0.0533961963580272
    "P10" = 0.03624253541791403
    "Q10" = 7.554140031384834
    "R10" = 45.324840188309
    "S10" = 0.0173411835292458
    "T10" = 0.01512434560348391
    "I11" = 0.3247643973172043
    "J11" = 0.4173092591090693
    "M11" = 79.17795566666666
    "N11" = 237.533867
    "O11" = 0.2937331182295834
    "P11" = 0.2990559347000084
    "Q11" = 41.55541514761878
    "R11" = 373.998736328569
    "S11" = 0.09539405911393375
    "T11" = 0.1247988105418307
    "G12" = 1.0751535
    "H12" = 2.150307
    "I12" = 0.6652969693706972
    "J12" = 0.5699199946567689
    "M12" = 63.46725166666666
    "N12" = 190.401755
    "O12" = 0.2354497988808272
    "P12" = 0.2397164477183668
    "Q12" = 68.23703776479749
    "R12" = 409.4222265887849
    "S12" = 0.1566440376343545
    "T12" = 0.1366191966027912
    "G13" = 1.0751535
    "H13" = 2.150307
    "I13" = 0.6652969693706972
    "J13" = 0.5699199946567689
    "O13" = 0.1779985000094065
    "P13" = 0.1812240584798697
    "Q13" = 51.58675193163749
    "R13" = 309.520511589825
    "S13" = 0.1184218626087881
    "T13" = 0.1032832144405253
    "G14" = 1.0751535
    "H14" = 2.150307
    "I14" = 0.6652969693706972
    "J14" = 0.5699199946567689
    "M14" = 64.53809233333334
    "N14" = 193.614277
    "O14" = 0.2394223865221556
    "P14" = 0.243761023683841
    "Q14" = 69.3883558555065
    "R14" = 416.330135133039
    "S14" = 0.1592869881526898
    "T14" = 0.1389242813154232
    "G15" = 1.0751535
    "H15" = 2.150307
    "I15" = 0.6652969693706972
    "J15" = 0.5699199946567689
    "M15" = 14.3933435
    "N15" = 28.786687
    "O15" = 0.0533961963580272
    "P15" = 0.03624253541791403
    "Q15" = 15.47505364072725
    "R15" = 61.90021456290899
    "S15" = 0.03552432761291816
    "T15" = 0.02065534559172532
    "G16" = 1.0751535
    "H16" = 2.150307
    "I16" = 0.6652969693706972
    "J16" = 0.5699199946567689
    "M16" = 79.17795566666666
    "N16" = 237.533867
    "O16" = 0.2937331182295834
    "P16" = 0.2990559347000084
    "Q16" = 85.12845615786148
    "R16" = 510.7707369471689
    "S16" = 0.1954197533619465
    "T16" = 0.1704379567063038
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Host "Applied $($updates.Count) cell updates."
